$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.782.51'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '2.378.35'
$ws.Range('E3').Value = '  -3.11%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '544.48'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.20'
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.541'
$ws.Range('E8').Value = '  -8.94%  '
$ws.Range('D9').Value = '2.375.93'
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.35'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.345'
$ws.Range('E13').Value = '  -1.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.50'
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '2.805.39'
$ws.Range('E15').Value = '  -3.00%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '60.627.05'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '2.376.74'
$ws.Range('E18').Value = '  -3.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.63'
$ws.Range('E19').Value = '  -3.58%  '
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '315.79'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('E24').Value = '  +2.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '62.76'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '2.492.40'
$ws.Range('E27').Value = '  -3.26%  '
$ws.Range('D28').Value = '0.0₃0928'
$ws.Range('E28').Value = '  -4.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.73'
$ws.Range('E29').Value = '  +2.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '521.08'
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('E31').Value = '  -3.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.99'
$ws.Range('E32').Value = '  -4.00%  '
$ws.Range('E33').Value = '  -3.67%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.48'
$ws.Range('E37').Value = '  -5.98%  '
$ws.Range('E38').Value = '  -3.34%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.01'
$ws.Range('E40').Value = '  -2.36%  '
$ws.Range('E41').Value = '  +1.55%  '
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.17'
$ws.Range('E43').Value = '  -5.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.28'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.22'
$ws.Range('E45').Value = '  -1.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '139.35'
$ws.Range('E46').Value = '  -4.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.55'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.28'
$ws.Range('E48').Value = '  -2.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0518'
$ws.Range('E49').Value = '  -1.89%  '
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0912'
$ws.Range('E51').Value = '  -2.48%  '
